$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 526
$ws1.Range("F5").Value = 130
$ws1.Range("F8").Value = 7129
$ws1.Range("F9").Value = 260
$ws1.Range("F11").Value = 3597
$ws1.Range("F14").Value = 266
$ws1.Range("F15").Value = 595

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 526
$ws4.Range("F6").Value = 130
$ws4.Range("F10").Value = 7129
$ws4.Range("F12").Value = 260
$ws4.Range("F14").Value = 3597
$ws4.Range("F17").Value = 266
$ws4.Range("F18").Value = 595
